$wb = $excel.ActiveWorkbook

# ---- Sheet "展览": remove the expired 2024-06-01 event (row 2), renumber, refresh counts ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows.Item(2).Delete()

$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(2,6).Value = 4208
$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(3,6).Value = 2402
$ws1.Cells.Item(4,1).Value = 3
$ws1.Cells.Item(4,6).Value = 478
$ws1.Cells.Item(5,1).Value = 4
$ws1.Cells.Item(5,6).Value = 18
$ws1.Cells.Item(6,1).Value = 5
$ws1.Cells.Item(6,6).Value = 43
$ws1.Cells.Item(7,1).Value = 6
$ws1.Cells.Item(7,6).Value = 44
$ws1.Cells.Item(8,1).Value = 7
$ws1.Cells.Item(8,6).Value = 212
$ws1.Cells.Item(9,1).Value = 8
$ws1.Cells.Item(9,6).Value = 122
$ws1.Cells.Item(10,1).Value = 9
$ws1.Cells.Item(10,6).Value = 117
$ws1.Cells.Item(11,1).Value = 10
$ws1.Cells.Item(11,6).Value = 148
$ws1.Cells.Item(12,1).Value = 11
$ws1.Cells.Item(12,6).Value = 1575
$ws1.Cells.Item(13,1).Value = 12
$ws1.Cells.Item(13,6).Value = 289
$ws1.Cells.Item(14,1).Value = 13
$ws1.Cells.Item(14,6).Value = 3208
$ws1.Cells.Item(15,1).Value = 14
$ws1.Cells.Item(15,6).Value = 215

# ---- Sheet "演出": refresh "想去人数" counts ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3,6).Value = 41
$ws2.Cells.Item(5,6).Value = 7

# ---- Sheet "全部类型": remove the expired 2024-06-01 event (row 2), renumber, refresh counts ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(2).Delete()

$ws4.Cells.Item(2,1).Value = 1
$ws4.Cells.Item(2,6).Value = 4208
$ws4.Cells.Item(3,1).Value = 2
$ws4.Cells.Item(3,6).Value = 2403
$ws4.Cells.Item(4,1).Value = 3
$ws4.Cells.Item(4,6).Value = 478
$ws4.Cells.Item(5,1).Value = 4
$ws4.Cells.Item(5,6).Value = 18
$ws4.Cells.Item(6,1).Value = 5
$ws4.Cells.Item(6,6).Value = 0
$ws4.Cells.Item(7,1).Value = 6
$ws4.Cells.Item(7,6).Value = 43
$ws4.Cells.Item(8,1).Value = 7
$ws4.Cells.Item(8,6).Value = 44
$ws4.Cells.Item(9,1).Value = 8
$ws4.Cells.Item(9,6).Value = 41
$ws4.Cells.Item(10,1).Value = 9
$ws4.Cells.Item(10,6).Value = 212
$ws4.Cells.Item(11,1).Value = 10
$ws4.Cells.Item(11,6).Value = 122
$ws4.Cells.Item(12,1).Value = 11
$ws4.Cells.Item(12,6).Value = 117
$ws4.Cells.Item(13,1).Value = 12
$ws4.Cells.Item(13,6).Value = 148
$ws4.Cells.Item(14,1).Value = 13
$ws4.Cells.Item(14,6).Value = 2
$ws4.Cells.Item(15,1).Value = 14
$ws4.Cells.Item(15,6).Value = 7
$ws4.Cells.Item(16,1).Value = 15
$ws4.Cells.Item(16,6).Value = 1575
$ws4.Cells.Item(17,1).Value = 16
$ws4.Cells.Item(17,6).Value = 289
$ws4.Cells.Item(18,1).Value = 17
$ws4.Cells.Item(18,6).Value = 3208
$ws4.Cells.Item(19,1).Value = 18
$ws4.Cells.Item(19,6).Value = 215
